# Insert two new rows at the top of the data table (rows 482-483), pushing
# all existing data rows down by two. This grows the table from
# A1:R606 to A1:R608, matching the "weekly" data refresh described in the
# commit message ("Fruta / hortaliza, semanal"): a new week's worth of
# observations (Primera / Segunda) is prepended, and the whole history
# shifts down accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift existing rows 482:606 down to 484:608, inserting two blank rows
# (with formatting carried over) at 482:483.
$ws.Rows("482:483").Insert()

# Populate the two newly inserted rows with the new week's data.
$ws.Cells.Item(482,1).Value = 6
$ws.Cells.Item(482,2).Value = 'Mercado Mayorista Lo Valledor de Santiago'
$ws.Cells.Item(482,3).Value = 'Metropolitana'
$ws.Cells.Item(482,4).Value = 44642
$ws.Cells.Item(482,5).Value = 13
$ws.Cells.Item(482,6).Value = 100112017
$ws.Cells.Item(482,7).Value = 'Apio'
$ws.Cells.Item(482,8).Value = 'Americana (o)'
$ws.Cells.Item(482,9).Value = 'Primera'
$ws.Cells.Item(482,10).Value = 2230
$ws.Cells.Item(482,11).Value = 7000
$ws.Cells.Item(482,12).Value = 8000
$ws.Cells.Item(482,13).Value = 7547
$ws.Cells.Item(482,14).Value = '$/docena de matas'
$ws.Cells.Item(482,15).Value = 'Región de Coquimbo'
$ws.Cells.Item(482,16).Value = 1258
$ws.Cells.Item(482,17).Value = 6
$ws.Cells.Item(482,18).Value = 'Hortaliza'

$ws.Cells.Item(483,1).Value = 6
$ws.Cells.Item(483,2).Value = 'Mercado Mayorista Lo Valledor de Santiago'
$ws.Cells.Item(483,3).Value = 'Metropolitana'
$ws.Cells.Item(483,4).Value = 44642
$ws.Cells.Item(483,5).Value = 13
$ws.Cells.Item(483,6).Value = 100112017
$ws.Cells.Item(483,7).Value = 'Apio'
$ws.Cells.Item(483,8).Value = 'Americana (o)'
$ws.Cells.Item(483,9).Value = 'Segunda'
$ws.Cells.Item(483,10).Value = 810
$ws.Cells.Item(483,11).Value = 5500
$ws.Cells.Item(483,12).Value = 6000
$ws.Cells.Item(483,13).Value = 5778
$ws.Cells.Item(483,14).Value = '$/docena de matas'
$ws.Cells.Item(483,15).Value = 'Región de Coquimbo'
$ws.Cells.Item(483,16).Value = 963
$ws.Cells.Item(483,17).Value = 6
$ws.Cells.Item(483,18).Value = 'Hortaliza'
